$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 4.2
$ws.Range("I5").Value = 5
$ws.Range("AE5").Value = 15

# Row 9
$ws.Range("K9").Value = 8.5
$ws.Range("L9").Value = 1.26
$ws.Range("M9").Value = 3.65
$ws.Range("N9").Value = 1.75
$ws.Range("O9").Value = 2.05
$ws.Range("P9").Value = 1.37
$ws.Range("Q9").Value = 2.92
$ws.Range("R9").Value = 1.7
$ws.Range("S9").Value = 2.05

# Row 11
$ws.Range("G11").Value = 1.75
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 1.04
$ws.Range("K11").Value = 9
$ws.Range("L11").Value = 1.23
$ws.Range("M11").Value = 3.9
$ws.Range("Q11").Value = 3
$ws.Range("X11").Value = 13
$ws.Range("AE11").Value = 15
$ws.Range("AF11").Value = 23
$ws.Range("AG11").Value = 15
$ws.Range("AI11").Value = 34

# Row 12
$ws.Range("G12").Value = 2.35
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 2.88
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 1.8
$ws.Range("Z12").Value = 9.5
$ws.Range("AA12").Value = 6
$ws.Range("AC12").Value = 51
$ws.Range("AJ12").Value = 34

# Row 16
$ws.Range("G16").Value = 1.07
$ws.Range("H16").Value = 7.5
$ws.Range("I16").Value = 25
$ws.Range("N16").Value = 1.26
$ws.Range("O16").Value = 3.5
$ws.Range("R16").Value = 2.44
$ws.Range("S16").Value = 1.49
$ws.Range("T16").Value = 10
$ws.Range("U16").Value = 6.3
$ws.Range("V16").Value = 11
$ws.Range("X16").Value = 9.5
$ws.Range("Y16").Value = 30
$ws.Range("Z16").Value = 23
$ws.Range("AA16").Value = 18.5
$ws.Range("AB16").Value = 35
$ws.Range("AC16").Value = 120
$ws.Range("AE16").Value = 90
$ws.Range("AF16").Value = 350
$ws.Range("AG16").Value = 80
$ws.Range("AI16").Value = 450
$ws.Range("AJ16").Value = 200

# Row 17
$ws.Range("G17").Value = 2.3
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 3.25
$ws.Range("K17").Value = 8.5
$ws.Range("L17").Value = 1.36
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 2.15
$ws.Range("O17").Value = 1.67
$ws.Range("R17").Value = 1.91
$ws.Range("S17").Value = 1.91
$ws.Range("T17").Value = 7
$ws.Range("X17").Value = 21
$ws.Range("Y17").Value = 34
$ws.Range("Z17").Value = 8.5
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 9
$ws.Range("AI17").Value = 29
$ws.Range("AJ17").Value = 41

# Row 20
$ws.Range("G20").Value = 3.25
$ws.Range("I20").Value = 2.1
$ws.Range("K20").Value = 9.5
$ws.Range("R20").Value = 1.8
$ws.Range("S20").Value = 1.91
$ws.Range("X20").Value = 29
$ws.Range("AC20").Value = 51
$ws.Range("AF20").Value = 10
$ws.Range("AG20").Value = 9
$ws.Range("AH20").Value = 19
$ws.Range("AI20").Value = 17

# Row 23
$ws.Range("N23").Value = 2.15
$ws.Range("O23").Value = 1.67

